$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 117, pushing the existing row 117 (and below)
# down by one. Excel's InsertRow copies the formatting of the row above into
# the newly inserted row, matching the style (s="2" date format) on column D.
$ws.Rows("117:117").Insert()

# New row 117 gets the new weekly record.
$ws.Range("A117").Value = 1
$ws.Range("B117").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C117").Value = "Arica y Parinacota"
$ws.Range("D117").Value = 44568
$ws.Range("E117").Value = 15
$ws.Range("F117").Value = "Fruta"
$ws.Range("G117").Value = 100108
$ws.Range("H117").Value = "Tropicales y subtropicales"
$ws.Range("I117").Value = 100108002
$ws.Range("J117").Value = "Mango"
$ws.Range("K117").Value = "Sin especificar"
$ws.Range("L117").Value = "Especial"
$ws.Range("M117").Value = 450
$ws.Range("N117").Value = 6000
$ws.Range("O117").Value = 6500
$ws.Range("P117").Value = 6250
$ws.Range("Q117").Value = "`$/bandeja 4 kilos"
$ws.Range("R117").Value = "Perú"
$ws.Range("S117").Value = 1562
$ws.Range("T117").Value = 4
